# Apply the RunControl / GlobalParams changes described by the commit
# "check AL issues in young and underfunded plans"

$wb = $excel.ActiveWorkbook

# --- RunControl sheet ---
$ws = $wb.Worksheets.Item("RunControl")

# Turn on "include" and "no_entrance" checks for the base DF100-1 run (row 6)
$ws.Range("C6").Value = $true
$ws.Range("K6").Value = $true

# Turn off "include" for the M1/M2/M3 (underfunded) runs, currently rows 17-19
$ws.Range("C17").Value = $false
$ws.Range("C18").Value = $false
$ws.Range("C19").Value = $false

# Remove the extra blank row (row 8) between the base run and the "O" block,
# shifting rows 9:19 up to become rows 8:18
$ws.Rows("8:8").Delete()

# Update the active selection to match the saved view
$ws.Range("D21").Select()

# --- GlobalParams sheet ---
$gp = $wb.Worksheets.Item("GlobalParams")
$gp.Range("A3").Value = 10
$gp.Range("A3").Select()

# Return focus to RunControl so it remains the active sheet/tab
$ws.Activate()
$ws.Range("D21").Select()
